# Atualizado por script em 31-10-2023 15:01
#
# The source site re-scraped the South Africa Premier League 2023-2024
# fixture list. A handful of already-recorded matches got reshuffled to
# different row positions (same match data, new row order) and four brand
# new matches were appended at the bottom of the sheet.
#
# Columns A (Indice) and E (data_partida) stay put for every existing row -
# only the match-specific payload in columns F:V (home team .. url_partida)
# moves between rows. New rows 70-73 are appended with the same column
# layout/styling as the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$FIRST_DATA_COL = 6   # column F
$LAST_DATA_COL  = 22  # column V

function Get-MatchRow([int]$row) {
    $vals = @()
    for ($c = $FIRST_DATA_COL; $c -le $LAST_DATA_COL; $c++) {
        $vals += ,$ws.Cells.Item($row, $c).Value()
    }
    return $vals
}

function Set-MatchRow([int]$row, $vals) {
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($row, $FIRST_DATA_COL + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# Re-shuffle the F:V (match) payload of the rows that moved. Row numbers,
# the "Indice" column (A) and the match date/time (E) are unchanged - only
# the home/away teams, scores, odds, timestamps and url move.
# ---------------------------------------------------------------------

# Rows 12 <-> 13 swap their match payload.
$orig12 = Get-MatchRow 12
$orig13 = Get-MatchRow 13
Set-MatchRow 12 $orig13
Set-MatchRow 13 $orig12

# Rows 15, 16, 17 rotate: 15<-17, 16<-15, 17<-16.
$orig15 = Get-MatchRow 15
$orig16 = Get-MatchRow 16
$orig17 = Get-MatchRow 17
Set-MatchRow 15 $orig17
Set-MatchRow 16 $orig15
Set-MatchRow 17 $orig16

# Rows 21 <-> 22 swap their match payload.
$orig21 = Get-MatchRow 21
$orig22 = Get-MatchRow 22
Set-MatchRow 21 $orig22
Set-MatchRow 22 $orig21

# Rows 37, 39, 41 rotate: 37<-41, 39<-37, 41<-39.
$orig37 = Get-MatchRow 37
$orig39 = Get-MatchRow 39
$orig41 = Get-MatchRow 41
Set-MatchRow 37 $orig41
Set-MatchRow 39 $orig37
Set-MatchRow 41 $orig39

# Rows 63, 64, 65 rotate: 63<-65, 64<-63, 65<-64.
$orig63 = Get-MatchRow 63
$orig64 = Get-MatchRow 64
$orig65 = Get-MatchRow 65
Set-MatchRow 63 $orig65
Set-MatchRow 64 $orig63
Set-MatchRow 65 $orig64

# ---------------------------------------------------------------------
# Append the four newly scraped matches as rows 70-73, matching the
# formatting (column A / E styles) already used by the last data row.
# ---------------------------------------------------------------------

$lastRow = 69
$newRowsCount = 4

$ws.Range("A$lastRow`:V$lastRow").Copy()
$ws.Range("A70:V$($lastRow + $newRowsCount)").PasteSpecial(-4122) # xlPasteFormats

$newRows = @(
    @{ Row = 70; A = 69; B = "south-africa"; C = "premier-league"; D = "2023-2024"; E = 45226.8125;
       F = "Orlando Pirates"; G = 1; H = "Polokwane"; I = 1;
       J = 1.88; K = "26/10/2023 12:31"; L = 1.5;  M = "27/10/2023 19:20";
       N = 3.17; O = "26/10/2023 12:31"; P = 4.04; Q = "27/10/2023 19:20";
       R = 4.41; S = "26/10/2023 12:31"; T = 7.39; U = "27/10/2023 19:20";
       V = "https://www.betexplorer.com/football/south-africa/premier-league/orlando-pirates-polokwane-city/E5eT3txj/" },

    @{ Row = 71; A = 70; B = "south-africa"; C = "premier-league"; D = "2023-2024"; E = 45227.64583333334;
       F = "Golden Arrows"; G = 2; H = "Kaizer Chiefs"; I = 1;
       J = 2.7;  K = "27/10/2023 15:12"; L = 3.35; M = "28/10/2023 15:25";
       N = 2.91; O = "27/10/2023 15:12"; P = 2.94; Q = "28/10/2023 15:25";
       R = 2.81; S = "27/10/2023 15:12"; T = 2.46; U = "28/10/2023 15:25";
       V = "https://www.betexplorer.com/football/south-africa/premier-league/golden-arrows-kaizer-chiefs/IZoJQ3q3/" },

    @{ Row = 72; A = 71; B = "south-africa"; C = "premier-league"; D = "2023-2024"; E = 45227.73958333334;
       F = "Chippa Utd."; G = 1; H = "AmaZulu"; I = 1;
       J = 2.24; K = "28/10/2023 15:13"; L = 2.29; M = "28/10/2023 15:16";
       N = 2.94; O = "28/10/2023 15:13"; P = 2.9;  Q = "28/10/2023 15:44";
       R = 3.52; S = "28/10/2023 15:13"; T = 3.56; U = "28/10/2023 15:16";
       V = "https://www.betexplorer.com/football/south-africa/premier-league/chippa-utd-amazulu/vqmFRqUd/" },

    @{ Row = 73; A = 72; B = "south-africa"; C = "premier-league"; D = "2023-2024"; E = 45228.60416666666;
       F = "Richards Bay"; G = 1; H = "Cape Town City"; I = 3;
       J = 2.86; K = "29/10/2023 07:59"; L = 2.96; M = "29/10/2023 14:24";
       N = 2.86; O = "29/10/2023 07:59"; P = 2.78; Q = "29/10/2023 14:24";
       R = 2.7;  S = "29/10/2023 07:59"; T = 2.9;  U = "29/10/2023 14:24";
       V = "https://www.betexplorer.com/football/south-africa/premier-league/richards-bay-cape-town-city/zBtOPNb9/" }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

foreach ($entry in $newRows) {
    $r = $entry.Row
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $entry[$col]
    }
}

Write-Host "Done. New dimension rows:" $ws.UsedRange.Rows.Count
